$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W8)
$ws1.Range("H2").Value = 14.14
$ws1.Range("L2").Value = 1.2

# Row 3 (W9)
$ws1.Range("H3").Value = 13.14
$ws1.Range("L3").Value = 0.8100000000000001

# Row 4 (W10)
$ws1.Range("H4").Value = 12.14
$ws1.Range("L4").Value = 0.82

# Row 5 (W11)
$ws1.Range("H5").Value = 10.42
$ws1.Range("L5").Value = 1.02

# Row 6 (W12)
$ws1.Range("H6").Value = 9.42
$ws1.Range("L6").Value = 1.07

# Row 7 (W13)
$ws1.Range("H7").Value = 8.880000000000001
$ws1.Range("L7").Value = 0.93

# Row 8 (W14)
$ws1.Range("H8").Value = 7.47
$ws1.Range("L8").Value = 0.86

# Row 9 (W15)
$ws1.Range("H9").Value = 6.47
$ws1.Range("L9").Value = 0.8100000000000001

# Row 10 (W16)
$ws1.Range("H10").Value = 5.3
$ws1.Range("L10").Value = 0.84

# Row 11 (W17)
$ws1.Range("H11").Value = 4.3
$ws1.Range("L11").Value = 1.18

# Row 12 (W18)
$ws1.Range("H12").Value = 3.3
$ws1.Range("L12").Value = 1.19

# Row 13 (W19)
$ws1.Range("H13").Value = 2.3
$ws1.Range("L13").Value = 0.9

# Row 14 (W20)
$ws1.Range("H14").Value = 1.3
$ws1.Range("L14").Value = 1.01

# Row 15 (W21)
$ws1.Range("H15").Value = 0.3
$ws1.Range("L15").Value = 1.15

# Row 17 (W23) - only Seasonality Index changes
$ws1.Range("L17").Value = 1

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "24"
$ws2.Range("B10").Style = "Normal"
